# Prasanna Skills Audit.xlsx update
# - Update "Date Reviewed:" label in B4 to include the review date
# - Add four new skill-audit rows (12-15) covering apex oracle, Slack,
#   Trello and One drive, reusing the formatting of the last existing
#   data row (row 17) for visual consistency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: record the review date -----------------------------------
$ws.Range("B4").Value = "Date Reviewed: 02/27/2020"

# --- Carry the formatting of the last populated row (17) down onto the
#     four still-blank rows (18-21) so the new rows look consistent with
#     the rest of the table (borders, fonts, etc.) ----------------------
$ws.Range("A17:H17").Copy() | Out-Null
$ws.Range("A18:H21").PasteSpecial(-4122) | Out-Null

# --- New data rows ------------------------------------------------------
# Row 18 -> item 12: apex oracle
$ws.Cells.Item(18, 1).Value = 12
$ws.Cells.Item(18, 2).Value = "apex oracle"
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = "Apex Oracle Online "
$ws.Cells.Item(18, 5).Value = "To manage and store data in a user friendly basis."
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(18, 7).Value = "29/04/2020"
$ws.Cells.Item(18, 8).Value = "Prasanna Shrestha"
$ws.Rows.Item(18).RowHeight = 31.2

# Row 19 -> item 13: Slack
$ws.Cells.Item(19, 1).Value = 13
$ws.Cells.Item(19, 2).Value = "Slack"
$ws.Cells.Item(19, 3).Value = 4
$ws.Cells.Item(19, 4).Value = "Slack Application"
$ws.Cells.Item(19, 5).Value = "For the communication of the team members"
$ws.Cells.Item(19, 6).Value = 5
$ws.Cells.Item(19, 7).Value = 43835
$ws.Cells.Item(19, 7).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(19, 8).Value = "Prasanna Shrestha"
$ws.Rows.Item(19).RowHeight = 31.2

# Row 20 -> item 14: Trello
$ws.Cells.Item(20, 1).Value = 14
$ws.Cells.Item(20, 2).Value = "Trello"
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = "Online course"
$ws.Cells.Item(20, 5).Value = "To make sure that the task is done in a proper organized way"
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(20, 7).Value = 43835
$ws.Cells.Item(20, 7).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(20, 8).Value = "Prasanna Shrestha"
$ws.Rows.Item(20).RowHeight = 31.2

# Row 21 -> item 15: One drive
$ws.Cells.Item(21, 1).Value = 15
$ws.Cells.Item(21, 2).Value = "One drive"
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(21, 4).Value = "Online course"
$ws.Cells.Item(21, 5).Value = "For the knowledge of files and folders storation"
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(21, 7).Value = 43835
$ws.Cells.Item(21, 7).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(21, 8).Value = "Prasanna Shrestha"
$ws.Rows.Item(21).RowHeight = 31.2

# Keep the active selection pointed at the header cell that was edited,
# matching how the workbook was left after the edit.
$ws.Range("B4").Select() | Out-Null
